$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'www.pineconedata.com'
$ws.Range("B3").Value = 'www.pineconedata.com'
$ws.Range("B4").Value = 'www.pineconedata.com'

$ws.Range("D2").Value = 'sampleCookie1'
$ws.Range("D3").Value = 'sampleCookie2'
$ws.Range("D4").Value = 'sampleCookie3'

$ws.Range("F2").Value = 'b"v10\xdcj\r\xe3m7\x17\x03b\xb7\xc7\xbe\x03\xf5t#\x03Z\n\xa8\xcd\xf6\xade\xf1\x05\x95\x0c\x88\xfc:\xd4\x91\x87|\xd7\xac\xdeG\x1c=\xbe\x97\x113\xa7}a\x12\x13\r\xde\xbbY\x82\x97n\xe9\xbd\xae\x8eW''\x02"'
$ws.Range("F3").Value = 'b''v10\xdcj\r\xe3m7\x17\x03b\xb7\xc7\xbe\x03\xf5t#\x03Z\n\xa8\xcd\xf6\xade\xf1\x05\x95\x0c\x88\xfc:\xd4$\x9a\xccsg-\xc4\xc6\\2y\xdf\x83D\x1afz=\x99\xbb2\xc4MVIgw\xe6\xc7\xb3\x1fS|\x9d\xed\xc2\xfdw\x8c/]@|\x92_j3S'''
$ws.Range("F4").Value = 'b''v10\xdcj\r\xe3m7\x17\x03b\xb7\xc7\xbe\x03\xf5t#\x03Z\n\xa8\xcd\xf6\xade\xf1\x05\x95\x0c\x88\xfc:\xd4\xbe\xfe\x98\x89\xc4x\x96\x96t\x12\x88^\x97F#E\xed\xa77\xf2\n\xd6\\2\x01\x9d\x98\x13\xa3\xd0t\x14'''

$ws.Range("G2").Value = '/'
$ws.Range("G3").Value = '/'
$ws.Range("G4").Value = '/'

$ws.Range("A2").Value = 13379971543596030
$ws.Range("A3").Value = 13379971543598250
$ws.Range("A4").Value = 13379971543600500

$ws.Range("H2").Value = 13380146637000000
$ws.Range("H3").Value = 13380146700000000
$ws.Range("H4").Value = 0

$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1

$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 1
$ws.Range("J4").Value = 1

$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0

$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 0

$ws.Range("M2").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("M4").Value = 0

$ws.Range("N2").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("N4").Value = 1

$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 2
$ws.Range("O4").Value = 0

$ws.Range("P2").Value = 2
$ws.Range("P3").Value = 2
$ws.Range("P4").Value = 2

$ws.Range("Q2").Value = 443
$ws.Range("Q3").Value = 443
$ws.Range("Q4").Value = 443

$ws.Range("R2").Value = 13379971543595960
$ws.Range("R3").Value = 13379971543598210
$ws.Range("R4").Value = 13379971543600460

$ws.Range("S2").Value = 3
$ws.Range("S3").Value = 3
$ws.Range("S4").Value = 3

$ws.Range("T2").Value = 1
$ws.Range("T3").Value = 1
$ws.Range("T4").Value = 1

